$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest crypto data refresh.
# D-column text values that look numeric need to be forced to Text format so Excel
# does not silently coerce them into numbers (matching the source inlineStr cells),
# then the style is reset back to Normal so no stray formatting is left behind.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '63.897.38'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.34%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.751.09'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +0.62%  '

$ws.Range("E4").Value = '  +0.02%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '575.04'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -0.59%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '157.67'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("E7").Value = '  +0.33%  '

$ws.Range("E8").Value = '  -1.50%  '

$ws.Range("E9").Value = '  -2.96%  '

$ws.Range("E10").Value = '  +1.81%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '5.66'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -15.79%  '

$ws.Range("E12").Value = '  -2.89%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '3.237.59'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.35%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '26.55'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -2.97%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '63.521.47'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -0.21%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.0000151'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -2.30%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '2.752.98'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -0.33%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '12.10'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +0.31%  '

$ws.Range("E19").Value = '  -2.22%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '355.40'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -1.73%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '6.70'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -3.97%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.32%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.534'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.26%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '65.17'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -2.31%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.171'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.41%  '

$ws.Range("E26").Value = '  -0.09%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '8.46'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -1.35%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0909'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -0.09%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.95'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -4.01%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '7.07'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -1.53%  '

$ws.Range("E31").Value = '  -1.22%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '168.45'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -4.15%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '20.15'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -2.15%  '

$ws.Range("E34").Value = '  +0.46%  '

$ws.Range("E35").Value = '  +0.21%  '

$ws.Range("E36").Value = '  +0.51%  '

$ws.Range("E37").Value = '  -1.69%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.983'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -2.50%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '6.17'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +6.26%  '

$ws.Range("E40").Value = '  -3.84%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '331.32'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -2.04%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '38.94'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -1.05%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '21.49'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -1.90%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.0587'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -2.21%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '21.50'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -2.85%  '

$ws.Range("E46").Value = '  -2.06%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.626'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -3.56%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '134.70'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -2.05%  '

$ws.Range("E49").Value = '  -0.86%  '

$ws.Range("E50").Value = '  +0.36%  '

$ws.Range("E51").Value = '  +0.07%  '
